$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.979.25'
$ws.Range("E2").Value = '  -2.73%  '
$ws.Range("D3").Value = '2.240.31'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'247.35"
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("E6").Value = '  -2.55%  '
$ws.Range("D7").Value = "'76.74"
$ws.Range("E7").Value = '  +3.58%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'0.621"
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").Value = "'40.95"
$ws.Range("E10").Value = '  +4.36%  '
$ws.Range("D11").Value = "'0.0950"
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("E12").Value = '  -3.30%  '
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").Value = '2.574.48'
$ws.Range("D15").Value = "'14.79"
$ws.Range("E15").Value = '  -3.30%  '
$ws.Range("D16").Value = "'0.861"
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("D17").Value = '2.246.74'
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("D18").Value = '41.906.25'
$ws.Range("E18").Value = '  -2.60%  '
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D20").Value = "'6.13"
$ws.Range("E20").Value = '  -2.93%  '
$ws.Range("D21").Value = "'71.72"
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("E22").Value = '  +4.30%  '
$ws.Range("D23").Value = "'230.54"
$ws.Range("E23").Value = '  -2.57%  '
$ws.Range("D24").Value = "'11.53"
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = "'3.67"
$ws.Range("E26").Value = '  -5.90%  '
$ws.Range("E27").Value = '  -5.15%  '
$ws.Range("D28").Value = "'7.19"
$ws.Range("E28").Value = '  +10.94%  '
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").Value = "'168.72"
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("D31").Value = "'20.55"
$ws.Range("E31").Value = '  -2.55%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = "'32.79"
$ws.Range("E32").Value = '  +3.84%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.0826"
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("E34").Value = '  -5.87%  '
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("D36").Value = "'4.52"
$ws.Range("E36").Value = '  -2.04%  '
$ws.Range("D37").Value = "'4.89"
$ws.Range("E37").Value = '  +2.57%  '
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("D39").Value = "'14.04"
$ws.Range("E39").Value = '  -3.35%  '
$ws.Range("E40").Value = '  -0.51%  '
$ws.Range("E41").Value = '  -6.95%  '
$ws.Range("D42").Value = "'112.82"
$ws.Range("E42").Value = '  +13.80%  '
$ws.Range("D43").Value = "'0.202"
$ws.Range("E43").Value = '  -7.37%  '
$ws.Range("D44").Value = "'61.25"
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("E45").Value = '  -4.13%  '
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("D47").Value = "'0.997"
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("D50").Value = "'4.31"
$ws.Range("E50").Value = '  -12.74%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = "'4.18"
$ws.Range("E51").Value = '  -2.17%  '
